# Sprint 1 Backlog Burndown - add initial estimates and actual hours
# (matches commit: "Updated The Burndown Chart / Added my initial estimates and actual hours")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8: "Create account/user page UI" - Initial Estimate / Week1 remaining / Week2 remaining
$ws.Range("D8").Value = 5
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 0

# Row 9: "Implement functionality for create account/user Page" - Initial Estimate
$ws.Range("D9").Value = 0

# Row 17: "Create Project Management Page UI ..." - Initial Estimate / Week1 remaining / Week2 remaining
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 4

# Row 18: "Implement functionality for Project Management Page ..." - Initial Estimate / Week1 remaining / Week2 remaining
$ws.Range("D18").Value = 10
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 13

# Move the selection/scroll to reflect the place the author was last working
$ws.Application.Goto($ws.Range("B8"))
$ws.Range("G18").Select()
